$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 data
$ws.Range("A12").Value = 1
$ws.Range("B12").Value = 9
$ws.Range("C12").Value = "Cauan De Oliveira Campos"
$ws.Range("D12").Value = 13
$ws.Range("E12").Value = "Ativo"
$ws.Range("F12").Value = 0.7
$ws.Range("G12").Value = 0.75
$ws.Range("H12").Value = 0.75
$ws.Range("I12").Value = 44960
$ws.Range("J12").Value = 45275
$ws.Range("K12").Value = 31026932
$ws.Range("L12").Value = 900086943667
$ws.Range("M12").Value = "000113555719 - 6"
$ws.Range("N12").Value = "8° ANO B INTEGRAL ANUAL"

# Row 13 data
$ws.Range("A13").Value = 1
$ws.Range("B13").Value = 39
$ws.Range("C13").Value = "Ian Lucas Alves Silva Moura"
$ws.Range("D13").Value = 14
$ws.Range("E13").Value = "Ativo"
$ws.Range("F13").Value = 0.73
$ws.Range("G13").Value = 0.85
$ws.Range("H13").Value = 0.82
$ws.Range("I13").Value = 45170
$ws.Range("J13").Value = 45275
$ws.Range("K13").Value = 30092959
$ws.Range("L13").Value = 900108744400
$ws.Range("M13").Value = "000112621744 - X"
$ws.Range("N13").Value = "8° ANO B INTEGRAL ANUAL"

# Apply the same date number format used in I2:J11 to the new date cells
$ws.Range("I12:J13").NumberFormat = $ws.Range("I11:J11").NumberFormat
